# Added support for date fields: "value" and "budget" rows on the
# "details" sheet, plus the matching shared-string entries and the
# column-width metadata that LibreOffice/Excel re-writes whenever a
# sheet's columns are touched.

$wb = $excel.ActiveWorkbook

# --- "toto" sheet: gains a <cols> default-width override ------------------
$toto = $wb.Worksheets.Item("toto")
$toto.Range($toto.Cells.Item(1, 1), $toto.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 7.666666667

# --- "tata" sheet: gains the same <cols> default-width override -----------
$tata = $wb.Worksheets.Item("tata")
$tata.Range($tata.Cells.Item(1, 1), $tata.Cells.Item(1, 1025)).EntireColumn.ColumnWidth = 7.666666667

# --- "details" sheet: two new rows (value / budget) ------------------------
$details = $wb.Worksheets.Item("details")

$details.Range("B15").Value = "value"
$details.Range("C15").Value = 50000

$details.Range("B16").Value = "budget"
$details.Range("C16").Value = 300

$details.Range("C16").Select()
